$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the greeting text for the 08:00-11:00 rule row (R10) from
# "Good Morning" to "GIT UPDATE".
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the active cell/selection recorded in the saved file.
$ws.Range("E8").Select()
